$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9552715654952076
$ws.Range("C2").Value = 0.7891373801916933

$ws.Range("B3").Value = 0.9680511182108626
$ws.Range("C3").Value = 0.8306709265175719

$ws.Range("B4").Value = 0.9744408945686901
$ws.Range("C4").Value = 0.7955271565495208

$ws.Range("B5").Value = 0.9648562300319489
$ws.Range("C5").Value = 0.8338658146964856

$ws.Range("B6").Value = 0.9712460063897763
$ws.Range("C6").Value = 0.8274760383386581
